$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns B:E to Text format while we overwrite values, so that
# numeric-looking strings (e.g. "314.82", "0.07286") are preserved exactly
# as text instead of being auto-coerced into floating point numbers.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "27.587.35"
$ws.Range("E2").Value = "  -2.53%  "

# Row 3
$ws.Range("D3").Value = "1.844.54"
$ws.Range("E3").Value = "  -1.51%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "314.82"
$ws.Range("E5").Value = "  -1.17%  "

# Row 6
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.05%  "

# Row 7
$ws.Range("D7").Value = "0.4274"
$ws.Range("E7").Value = "  -2.94%  "

# Row 8
$ws.Range("D8").Value = "0.3653"
$ws.Range("E8").Value = "  -1.33%  "

# Row 9
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07286"
$ws.Range("E9").Value = "  -3.07%  "

# Row 10
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "0.8962"
$ws.Range("E10").Value = "  -4.64%  "

# Row 11
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "20.69"
$ws.Range("E11").Value = "  -3.46%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.840.40"
$ws.Range("E12").Value = "  -3.80%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.394"
$ws.Range("E13").Value = "  -1.20%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "6.578"
$ws.Range("E14").Value = "  -2.29%  "

# Row 15
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "0.06869"
$ws.Range("E15").Value = "  +0.19%  "

# Row 16
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.01%  "

# Row 17
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "78.58"
$ws.Range("E17").Value = "  -4.47%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.000008886"
$ws.Range("E18").Value = "  -1.83%  "

# Row 19
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.07%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "15.55"
$ws.Range("E20").Value = "  -2.48%  "

# Row 21
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "27.607.31"
$ws.Range("E21").Value = "  -2.45%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "4.993"
$ws.Range("E22").Value = "  -2.81%  "

# Row 23
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "10.54"
$ws.Range("E23").Value = "  -3.13%  "

# Row 24
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "2.060.87"
$ws.Range("E24").Value = "  -3.10%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.039"
$ws.Range("E25").Value = "  +0.44%  "

# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "155.05"
$ws.Range("E26").Value = "  +0.06%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "18.46"
$ws.Range("E27").Value = "  +0.44%  "

# Row 28
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "5.244"
$ws.Range("E28").Value = "  -1.81%  "

# Row 29
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "115.36"
$ws.Range("E29").Value = "  +1.39%  "

# Row 30
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "1.839"
$ws.Range("E30").Value = "  +6.13%  "

# Row 31
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.08898"
$ws.Range("E31").Value = "  -1.67%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "0.7870"
$ws.Range("E32").Value = "  -1.75%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "4.579"
$ws.Range("E33").Value = "  -5.61%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.961"
$ws.Range("E34").Value = "  +1.68%  "

# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.119"
$ws.Range("E35").Value = "  -4.98%  "

# Row 36
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "0.9999"
$ws.Range("E36").Value = "  -0.11%  "

# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.05454"
$ws.Range("E37").Value = "  +0.16%  "

# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.102"
$ws.Range("E38").Value = "  -2.11%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01930"
$ws.Range("E39").Value = "  -1.50%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.775"
$ws.Range("E40").Value = "  -5.02%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.5074"
$ws.Range("E41").Value = "  -3.64%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "6.832"
$ws.Range("E42").Value = "  -4.19%  "

# Row 43
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.1650"
$ws.Range("E43").Value = "  -2.16%  "

# Row 44
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "8.264"
$ws.Range("E44").Value = "  -5.84%  "

# Row 45
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "0.06637"
$ws.Range("E45").Value = "  -1.76%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "10.42"
$ws.Range("E46").Value = "  -1.56%  "

# Row 47
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.4721"
$ws.Range("E47").Value = "  -3.46%  "

# Row 48
$ws.Range("E48").Value = "  -2.03%  "

# Row 49
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "0.9998"
$ws.Range("E49").Value = "  -0.05%  "

# Row 50
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.869"
$ws.Range("E50").Value = "  -6.08%  "

# Row 51
$ws.Range("D51").Value = "1.636"
$ws.Range("E51").Value = "  -2.82%  "

# Restore the default cell style so the written cells match the original
# (unstyled) formatting instead of keeping the temporary Text number format.
$dataRange.Style = "Normal"

